$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: Add new header columns V1/W1 (copy format from U1, keep the same shared style s="1") ----
$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("U1").Copy()
$ws.Range("W1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("V1").Value = "Posesión Local ().2"
$ws.Range("W1").Value = "Posesión Visita ().2"
$excel.CutCopyMode = 0

# ---- Step 2: Fix up goals-by-half (1T/2T) values on existing rows (stat corrections) ----
$ws.Cells.Item(122, 13).Value = 1   # M122
$ws.Cells.Item(122, 15).Value = 1   # O122
$ws.Cells.Item(127, 13).Value = 1   # M127
$ws.Cells.Item(127, 14).Value = 1   # N127
$ws.Cells.Item(127, 15).Value = 1   # O127
$ws.Cells.Item(127, 16).Value = 1   # P127
$ws.Cells.Item(128, 14).Value = 2   # N128
$ws.Cells.Item(128, 16).Value = 0   # P128
$ws.Cells.Item(129, 13).Value = 1   # M129
$ws.Cells.Item(129, 15).Value = 0   # O129
$ws.Cells.Item(132, 13).Value = 2   # M132
$ws.Cells.Item(132, 15).Value = 0   # O132
$ws.Cells.Item(133, 13).Value = 1   # M133
$ws.Cells.Item(133, 15).Value = 0   # O133
$ws.Cells.Item(134, 13).Value = 1   # M134
$ws.Cells.Item(134, 15).Value = 3   # O134
$ws.Cells.Item(135, 13).Value = 1   # M135
$ws.Cells.Item(135, 15).Value = 0   # O135
$ws.Cells.Item(136, 13).Value = 1   # M136
$ws.Cells.Item(136, 14).Value = 1   # N136
$ws.Cells.Item(136, 15).Value = 2   # O136
$ws.Cells.Item(136, 16).Value = 1   # P136
$ws.Cells.Item(139, 13).Value = 1   # M139
$ws.Cells.Item(139, 14).Value = 1   # N139
$ws.Cells.Item(139, 15).Value = 0   # O139
$ws.Cells.Item(139, 16).Value = 1   # P139
$ws.Cells.Item(140, 14).Value = 2   # N140
$ws.Cells.Item(140, 16).Value = 0   # P140
$ws.Cells.Item(142, 14).Value = 1   # N142
$ws.Cells.Item(142, 16).Value = 1   # P142
$ws.Cells.Item(145, 13).Value = 1   # M145
$ws.Cells.Item(145, 15).Value = 0   # O145
$ws.Cells.Item(146, 13).Value = 1   # M146
$ws.Cells.Item(146, 15).Value = 1   # O146
$ws.Cells.Item(147, 13).Value = 2   # M147
$ws.Cells.Item(147, 14).Value = 1   # N147
$ws.Cells.Item(147, 15).Value = 1   # O147
$ws.Cells.Item(147, 16).Value = 0   # P147
$ws.Cells.Item(148, 13).Value = 1   # M148
$ws.Cells.Item(148, 14).Value = 1   # N148
$ws.Cells.Item(148, 15).Value = 0   # O148
$ws.Cells.Item(148, 16).Value = 0   # P148
$ws.Cells.Item(149, 13).Value = 1   # M149
$ws.Cells.Item(149, 15).Value = 0   # O149
$ws.Cells.Item(150, 13).Value = 1   # M150
$ws.Cells.Item(150, 15).Value = 2   # O150
$ws.Cells.Item(151, 13).Value = 1   # M151
$ws.Cells.Item(151, 14).Value = 1   # N151
$ws.Cells.Item(151, 15).Value = 0   # O151
$ws.Cells.Item(151, 16).Value = 1   # P151
$ws.Cells.Item(152, 13).Value = 1   # M152
$ws.Cells.Item(152, 15).Value = 2   # O152
$ws.Cells.Item(153, 14).Value = 1   # N153
$ws.Cells.Item(153, 16).Value = 0   # P153
$ws.Cells.Item(155, 13).Value = 2   # M155
$ws.Cells.Item(155, 14).Value = 1   # N155
$ws.Cells.Item(155, 15).Value = 0   # O155
$ws.Cells.Item(155, 16).Value = 0   # P155
$ws.Cells.Item(158, 13).Value = 1   # M158
$ws.Cells.Item(158, 15).Value = 0   # O158
$ws.Cells.Item(159, 14).Value = 2   # N159
$ws.Cells.Item(159, 16).Value = 0   # P159
$ws.Cells.Item(160, 14).Value = 1   # N160
$ws.Cells.Item(160, 16).Value = 0   # P160
$ws.Cells.Item(161, 14).Value = 1   # N161
$ws.Cells.Item(161, 16).Value = 0   # P161
$ws.Cells.Item(162, 13).Value = 1   # M162
$ws.Cells.Item(162, 15).Value = 1   # O162
$ws.Cells.Item(163, 13).Value = 1   # M163
$ws.Cells.Item(163, 15).Value = 1   # O163
$ws.Cells.Item(164, 14).Value = 1   # N164
$ws.Cells.Item(164, 16).Value = 1   # P164
$ws.Cells.Item(165, 13).Value = 2   # M165
$ws.Cells.Item(165, 14).Value = 1   # N165
$ws.Cells.Item(165, 15).Value = 1   # O165
$ws.Cells.Item(165, 16).Value = 0   # P165

# ---- Step 3: Append the new fixtures (rows 166-175) ----
# Row 166: Sao Paulo vs Vitoria
$ws.Cells.Item(166, 1).NumberFormat = "@"
$ws.Cells.Item(166, 1).Value = "2025-08-09"
$ws.Cells.Item(166, 1).Style = "Normal"
$ws.Cells.Item(166, 2).Value = "Sao Paulo"
$ws.Cells.Item(166, 3).Value = "Vitoria"
$ws.Cells.Item(166, 4).Value = 2
$ws.Cells.Item(166, 5).Value = 0
$ws.Cells.Item(166, 6).Value = 1351226
$ws.Cells.Item(166, 7).Value = 4
$ws.Cells.Item(166, 8).Value = 3
$ws.Cells.Item(166, 9).Value = 3
$ws.Cells.Item(166, 10).Value = 3
$ws.Cells.Item(166, 11).Value = 0
$ws.Cells.Item(166, 12).Value = 0
$ws.Cells.Item(166, 13).Value = 1
$ws.Cells.Item(166, 14).Value = 0
$ws.Cells.Item(166, 15).Value = 1
$ws.Cells.Item(166, 16).Value = 0
$ws.Cells.Item(166, 17).Value = 56
$ws.Cells.Item(166, 18).Value = 44
$ws.Cells.Item(166, 19).Value = "L"

# Row 167: Flamengo vs Mirassol
$ws.Cells.Item(167, 1).NumberFormat = "@"
$ws.Cells.Item(167, 1).Value = "2025-08-09"
$ws.Cells.Item(167, 1).Style = "Normal"
$ws.Cells.Item(167, 2).Value = "Flamengo"
$ws.Cells.Item(167, 3).Value = "Mirassol"
$ws.Cells.Item(167, 4).Value = 2
$ws.Cells.Item(167, 5).Value = 1
$ws.Cells.Item(167, 6).Value = 1351223
$ws.Cells.Item(167, 7).Value = 11
$ws.Cells.Item(167, 8).Value = 10
$ws.Cells.Item(167, 9).Value = 1
$ws.Cells.Item(167, 10).Value = 1
$ws.Cells.Item(167, 11).Value = 0
$ws.Cells.Item(167, 12).Value = 0
$ws.Cells.Item(167, 13).Value = 1
$ws.Cells.Item(167, 14).Value = 0
$ws.Cells.Item(167, 15).Value = 1
$ws.Cells.Item(167, 16).Value = 1
$ws.Cells.Item(167, 17).Value = 49
$ws.Cells.Item(167, 18).Value = 51
$ws.Cells.Item(167, 19).Value = "L"

# Row 168: RB Bragantino vs Internacional
$ws.Cells.Item(168, 1).NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = "2025-08-09"
$ws.Cells.Item(168, 1).Style = "Normal"
$ws.Cells.Item(168, 2).Value = "RB Bragantino"
$ws.Cells.Item(168, 3).Value = "Internacional"
$ws.Cells.Item(168, 4).Value = 1
$ws.Cells.Item(168, 5).Value = 3
$ws.Cells.Item(168, 6).Value = 1351227
$ws.Cells.Item(168, 7).Value = 8
$ws.Cells.Item(168, 8).Value = 6
$ws.Cells.Item(168, 9).Value = 1
$ws.Cells.Item(168, 10).Value = 2
$ws.Cells.Item(168, 11).Value = 0
$ws.Cells.Item(168, 12).Value = 0
$ws.Cells.Item(168, 13).Value = 0
$ws.Cells.Item(168, 14).Value = 2
$ws.Cells.Item(168, 15).Value = 1
$ws.Cells.Item(168, 16).Value = 1
$ws.Cells.Item(168, 17).Value = 50
$ws.Cells.Item(168, 18).Value = 50
$ws.Cells.Item(168, 19).Value = "V"

# Row 169: Fortaleza EC vs Botafogo
$ws.Cells.Item(169, 1).NumberFormat = "@"
$ws.Cells.Item(169, 1).Value = "2025-08-09"
$ws.Cells.Item(169, 1).Style = "Normal"
$ws.Cells.Item(169, 2).Value = "Fortaleza EC"
$ws.Cells.Item(169, 3).Value = "Botafogo"
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = 5
$ws.Cells.Item(169, 6).Value = 1351231
$ws.Cells.Item(169, 7).Value = 2
$ws.Cells.Item(169, 8).Value = 8
$ws.Cells.Item(169, 9).Value = 2
$ws.Cells.Item(169, 10).Value = 0
$ws.Cells.Item(169, 11).Value = 1
$ws.Cells.Item(169, 12).Value = 0
$ws.Cells.Item(169, 13).Value = 0
$ws.Cells.Item(169, 14).Value = 2
$ws.Cells.Item(169, 15).Value = 0
$ws.Cells.Item(169, 16).Value = 3
$ws.Cells.Item(169, 17).Value = 43
$ws.Cells.Item(169, 18).Value = 57
$ws.Cells.Item(169, 19).Value = "V"

# Row 170: Bahia vs Fluminense
$ws.Cells.Item(170, 1).NumberFormat = "@"
$ws.Cells.Item(170, 1).Value = "2025-08-10"
$ws.Cells.Item(170, 1).Style = "Normal"
$ws.Cells.Item(170, 2).Value = "Bahia"
$ws.Cells.Item(170, 3).Value = "Fluminense"
$ws.Cells.Item(170, 4).Value = 3
$ws.Cells.Item(170, 5).Value = 3
$ws.Cells.Item(170, 6).Value = 1351230
$ws.Cells.Item(170, 7).Value = 8
$ws.Cells.Item(170, 8).Value = 2
$ws.Cells.Item(170, 9).Value = 2
$ws.Cells.Item(170, 10).Value = 2
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(170, 12).Value = 1
$ws.Cells.Item(170, 13).Value = 2
$ws.Cells.Item(170, 14).Value = 1
$ws.Cells.Item(170, 15).Value = 1
$ws.Cells.Item(170, 16).Value = 2
$ws.Cells.Item(170, 17).Value = 60
$ws.Cells.Item(170, 18).Value = 40
$ws.Cells.Item(170, 19).Value = "E"

# Row 171: Palmeiras vs Ceara
$ws.Cells.Item(171, 1).NumberFormat = "@"
$ws.Cells.Item(171, 1).Value = "2025-08-10"
$ws.Cells.Item(171, 1).Style = "Normal"
$ws.Cells.Item(171, 2).Value = "Palmeiras"
$ws.Cells.Item(171, 3).Value = "Ceara"
$ws.Cells.Item(171, 4).Value = 2
$ws.Cells.Item(171, 5).Value = 1
$ws.Cells.Item(171, 6).Value = 1351225
$ws.Cells.Item(171, 7).Value = 3
$ws.Cells.Item(171, 8).Value = 3
$ws.Cells.Item(171, 9).Value = 1
$ws.Cells.Item(171, 10).Value = 3
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 0
$ws.Cells.Item(171, 13).Value = 0
$ws.Cells.Item(171, 14).Value = 0
$ws.Cells.Item(171, 15).Value = 2
$ws.Cells.Item(171, 16).Value = 1
$ws.Cells.Item(171, 17).Value = 62
$ws.Cells.Item(171, 18).Value = 38
$ws.Cells.Item(171, 19).Value = "L"

# Row 172: Vasco DA Gama vs Atletico-MG
$ws.Cells.Item(172, 1).NumberFormat = "@"
$ws.Cells.Item(172, 1).Value = "2025-08-10"
$ws.Cells.Item(172, 1).Style = "Normal"
$ws.Cells.Item(172, 2).Value = "Vasco DA Gama"
$ws.Cells.Item(172, 3).Value = "Atletico-MG"
$ws.Cells.Item(172, 4).Value = 1
$ws.Cells.Item(172, 5).Value = 1
$ws.Cells.Item(172, 6).Value = 1351224
$ws.Cells.Item(172, 7).Value = 6
$ws.Cells.Item(172, 8).Value = 10
$ws.Cells.Item(172, 9).Value = 3
$ws.Cells.Item(172, 10).Value = 1
$ws.Cells.Item(172, 11).Value = 0
$ws.Cells.Item(172, 12).Value = 0
$ws.Cells.Item(172, 13).Value = 1
$ws.Cells.Item(172, 14).Value = 1
$ws.Cells.Item(172, 15).Value = 0
$ws.Cells.Item(172, 16).Value = 0
$ws.Cells.Item(172, 17).Value = 55
$ws.Cells.Item(172, 18).Value = 45
$ws.Cells.Item(172, 19).Value = "E"

# Row 173: Cruzeiro vs Santos
$ws.Cells.Item(173, 1).NumberFormat = "@"
$ws.Cells.Item(173, 1).Value = "2025-08-10"
$ws.Cells.Item(173, 1).Style = "Normal"
$ws.Cells.Item(173, 2).Value = "Cruzeiro"
$ws.Cells.Item(173, 3).Value = "Santos"
$ws.Cells.Item(173, 4).Value = 1
$ws.Cells.Item(173, 5).Value = 2
$ws.Cells.Item(173, 6).Value = 1351228
$ws.Cells.Item(173, 7).Value = 15
$ws.Cells.Item(173, 8).Value = 3
$ws.Cells.Item(173, 9).Value = 5
$ws.Cells.Item(173, 10).Value = 4
$ws.Cells.Item(173, 11).Value = 0
$ws.Cells.Item(173, 12).Value = 0
$ws.Cells.Item(173, 13).Value = 1
$ws.Cells.Item(173, 14).Value = 0
$ws.Cells.Item(173, 15).Value = 0
$ws.Cells.Item(173, 16).Value = 2
$ws.Cells.Item(173, 17).Value = 57
$ws.Cells.Item(173, 18).Value = 43
$ws.Cells.Item(173, 19).Value = "V"

# Row 174: Gremio vs Sport Recife
$ws.Cells.Item(174, 1).NumberFormat = "@"
$ws.Cells.Item(174, 1).Value = "2025-08-10"
$ws.Cells.Item(174, 1).Style = "Normal"
$ws.Cells.Item(174, 2).Value = "Gremio"
$ws.Cells.Item(174, 3).Value = "Sport Recife"
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 1
$ws.Cells.Item(174, 6).Value = 1351229
$ws.Cells.Item(174, 7).Value = 3
$ws.Cells.Item(174, 8).Value = 2
$ws.Cells.Item(174, 9).Value = 2
$ws.Cells.Item(174, 10).Value = 0
$ws.Cells.Item(174, 11).Value = 0
$ws.Cells.Item(174, 12).Value = 0
$ws.Cells.Item(174, 13).Value = 0
$ws.Cells.Item(174, 14).Value = 0
$ws.Cells.Item(174, 15).Value = 0
$ws.Cells.Item(174, 16).Value = 1
$ws.Cells.Item(174, 17).Value = 59
$ws.Cells.Item(174, 18).Value = 41
$ws.Cells.Item(174, 19).Value = "V"

# Row 175: Juventude vs Corinthians
$ws.Cells.Item(175, 1).NumberFormat = "@"
$ws.Cells.Item(175, 1).Value = "2025-08-11"
$ws.Cells.Item(175, 1).Style = "Normal"
$ws.Cells.Item(175, 2).Value = "Juventude"
$ws.Cells.Item(175, 3).Value = "Corinthians"
$ws.Cells.Item(175, 4).Value = 2
$ws.Cells.Item(175, 5).Value = 1
$ws.Cells.Item(175, 6).Value = 1351232
$ws.Cells.Item(175, 7).Value = 3
$ws.Cells.Item(175, 8).Value = 2
$ws.Cells.Item(175, 9).Value = 6
$ws.Cells.Item(175, 10).Value = 5
$ws.Cells.Item(175, 11).Value = 0
$ws.Cells.Item(175, 12).Value = 1
$ws.Cells.Item(175, 13).Value = 1
$ws.Cells.Item(175, 14).Value = 0
$ws.Cells.Item(175, 15).Value = 1
$ws.Cells.Item(175, 16).Value = 1
$ws.Cells.Item(175, 17).Value = 37
$ws.Cells.Item(175, 18).Value = 63
$ws.Cells.Item(175, 19).Value = "L"

